$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '309.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-3.24%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '48.98'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.10%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.97%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07771'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-4.11%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.508'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.12%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.390'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '16.38%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.559'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-6.22%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-6.42%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1994'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.96%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.04684'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.78%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09408'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.29%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1048'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.14%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-4.64%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04187'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-2.62%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005798'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.11%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2,018.73%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.330'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.03%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.240'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-8.05%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3501'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.20%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.928'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-3.03%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1340'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-5.72%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.77%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001272'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.94%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004052'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-4.55%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001352'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.00%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02595'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-3.24%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05874'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '5.33%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01095'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '73.64%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007914'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.81%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1422'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.10%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008455'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '9.76%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008354'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '3.12%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3101'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-2.91%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007054'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.82%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000752'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.12%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05302'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-13.32%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002626'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-34.44%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002105'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.12%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002005'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.12%'
